# Update cryptos list data (prices + 1h volume change) per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.062.54"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "1.650.46"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").Value = "'218.07"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").Value = "'0.5292"
$ws.Range("E6").Value = "  +1.67%  "
$ws.Range("D7").Value = "'1.003"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("D9").Value = "'0.06304"
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").Value = "'20.32"
$ws.Range("E10").Value = "  -3.33%  "
$ws.Range("D11").Value = "'0.07736"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.679.48"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.475"
$ws.Range("E13").Value = "  +1.24%  "
$ws.Range("D14").Value = "'0.5449"
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").Value = "0.0₅8106"
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("D16").Value = "'65.06"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").Value = "26.061.96"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("E19").Value = "  -2.34%  "
$ws.Range("D20").Value = "'193.17"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").Value = "'10.03"
$ws.Range("E21").Value = "  -1.02%  "
$ws.Range("D22").Value = "'5.991"
$ws.Range("E22").Value = "  -1.51%  "
$ws.Range("D23").Value = "'1.004"
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").Value = "'140.09"
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("D25").Value = "'0.1239"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").Value = "'7.231"
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("D27").Value = "'16.17"
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("D28").Value = "'1.434"
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("D29").Value = "'0.05914"
$ws.Range("E29").Value = "  -1.44%  "
$ws.Range("D30").Value = "'1.276"
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("D31").Value = "'3.499"
$ws.Range("E31").Value = "  -1.83%  "
$ws.Range("D32").Value = "'3.238"
$ws.Range("E32").Value = "  -2.65%  "
$ws.Range("D33").Value = "'1.543"
$ws.Range("E33").Value = "  -6.08%  "
$ws.Range("D34").Value = "'2.413"
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("D35").Value = "'0.9411"
$ws.Range("E35").Value = "  -3.91%  "
$ws.Range("D36").Value = "'2.756"
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("D37").Value = "'0.5688"
$ws.Range("E37").Value = "  -3.21%  "
$ws.Range("D38").Value = "'0.01605"
$ws.Range("E38").Value = "  +1.23%  "
$ws.Range("E39").Value = "  -1.51%  "
$ws.Range("D40").Value = "'0.8437"
$ws.Range("E40").Value = "  -2.37%  "
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").Value = "1.007.87"
$ws.Range("E42").Value = "  -2.72%  "
$ws.Range("E43").Value = "  +1.15%  "
$ws.Range("D44").Value = "1.799.90"
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").Value = "'56.85"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("D47").Value = "'1.007"
$ws.Range("E47").Value = "  +0.37%  "
$ws.Range("D48").Value = "'0.4295"
$ws.Range("E48").Value = "  +1.61%  "
$ws.Range("D49").Value = "'1.481"
$ws.Range("E49").Value = "  +2.23%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.05152"
$ws.Range("E50").Value = "  -0.58%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.792"
$ws.Range("E51").Value = "  -3.38%  "
